# Append the new daily COVID-19 data row (11/6/2020, serial 43993) to the
# bottom of the "Tabela1" table on the first worksheet, mirroring a user
# who typed a new row right under the table in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item("Tabela1")

# Grow the table by one row - this also extends the table ref / autofilter
# range and the worksheet's used-range dimension automatically.
$newRow = $lo.ListRows.Add()

# Pick up the same look (number formats / borders / font) as the previous
# data row, then overwrite the new row with numeric values.
$ws.Range("A89:J89").Copy() | Out-Null
$newRow.Range.PasteSpecial(-4122) | Out-Null

$newRow.Range.Item(1).Value = 43993
$newRow.Range.Item(2).Value = 86328
$newRow.Range.Item(3).Value = 702
$newRow.Range.Item(4).Value = 1490
$newRow.Range.Item(5).Value = 2
$newRow.Range.Item(6).Value = 6
$newRow.Range.Item(7).Value = 0
$newRow.Range.Item(8).Value = 0
$newRow.Range.Item(9).Value = 109
$newRow.Range.Item(10).Value = 0

$newRow.Range.Select() | Out-Null
